$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 data (new fixture values)
$ws.Range("C5").Value = "eogaxjynuzvkhno@gmail.com"
$ws.Range("D5").Value = "qfgqfXQUGL5"
$ws.Range("F5").ClearContents()

# Remove rows 6 and 7 entirely (shrinks used range / dimension to A1:F5)
$ws.Range("A6:F7").EntireRow.Delete()

# Move the active selection to F13, matching the saved view state
$ws.Range("F13").Select()
